$wb = $excel.ActiveWorkbook

# --- Sheet "Forecast Comparison": D3 numeric value 12 -> 13 ---
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D3").Value = 13

# --- Sheet "Summary": B9, B10, B11 are text-typed numeric strings ---
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "307"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "141"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "58"
